$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2183
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 0.1216690927753454
$ws.Range("E2").Value = 0.005599832721877925
$ws.Range("B3").Value = 2205
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0.0443432380026855
$ws.Range("E3").Value = 0.003034898418434704
$ws.Range("B4").Value = 2232
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0.09183948041187653
$ws.Range("E4").Value = 0.004683980786033763
$ws.Range("B5").Value = 2542.166758096668
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.01383495386696538
$ws.Range("E5").Value = 0.0008173620906235424
$ws.Range("B6").Value = 2551.502019868145
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.01373841949095679
$ws.Range("E6").Value = 0.0006814880343547023
$ws.Range("B7").Value = 2564.932497649872
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0.01819099391571898
$ws.Range("E7").Value = 0.0006972695510192858
$ws.Range("B8").Value = 2581.472193600247
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.01740303898101394
$ws.Range("E8").Value = 0.0006407842842457109
$ws.Range("B9").Value = 2610.832823110164
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 0.01810023613242792
$ws.Range("E9").Value = 0.004188484394280823
$ws.Range("B10").Value = 2622.008751622352
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0.02696907165021489
$ws.Range("E10").Value = 0.001026990540070901
$ws.Range("B11").Value = 2634.905719285064
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0.003117341027255372
$ws.Range("E11").Value = 0.0006073936095801519
$ws.Range("B12").Value = 2644.795530877373
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 0.01746763940385893
$ws.Range("E12").Value = 0.003700771060139565
$ws.Range("B13").Value = 2659.460038849562
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 0.02298617066354507
$ws.Range("E13").Value = 0.003978375691767457
$ws.Range("B14").Value = 2675.30036084316
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 0.01570416697021396
$ws.Range("E14").Value = 0.003962733721455892
$ws.Range("B15").Value = 2705.206271534023
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.01812844785914926
$ws.Range("E15").Value = 0.0008985371453171434
$ws.Range("B16").Value = 2715.278795007966
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 0.01534431112197311
$ws.Range("E16").Value = 0.0006147424095827502
$ws.Range("B17").Value = 2729.747422470061
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0.01047952597234795
$ws.Range("E17").Value = 0.0005900059890142009
$ws.Range("B18").Value = 2751.130289051061
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 0.04287844891871738
$ws.Range("E18").Value = 0.007457121551081283
$ws.Range("B19").Value = 2757.680085527932
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.01406444463080108
$ws.Range("E19").Value = 0.001024149975157255
$ws.Range("B20").Value = 2777.181176406648
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 0.01955736703752672
$ws.Range("E20").Value = 0.0007471648265006857
$ws.Range("B21").Value = 2787
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 0.04851787932698472
$ws.Range("E21").Value = 0.00121145340602076
$ws.Range("B22").Value = 2827.101259227846
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 0.003286920456302312
$ws.Range("E22").Value = 0.0003118027310326304
$ws.Range("B23").Value = 2841.959721907576
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 0.001179037588699636
$ws.Range("E23").Value = 0.0002591291403735464
$ws.Range("B24").Value = 2869.514732511383
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0.01742527498144795
$ws.Range("E24").Value = 0.0006691873879849242
$ws.Range("B25").Value = 2887.503509710963
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 0.01832457211463998
$ws.Range("E25").Value = 0.001638985953500047
$ws.Range("B26").Value = 2894.579211705925
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.0186799791228335
$ws.Range("E26").Value = 0.002949470387815886
$ws.Range("B27").Value = 2907.129675101117
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 0.01174384824347041
$ws.Range("E27").Value = 0.0008170322442685769
$ws.Range("B28").Value = 2924.234342911453
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0.02442874112444123
$ws.Range("E28").Value = 0.002682925325805757
$ws.Range("B29").Value = 2936.124980912122
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 0.02124869516641243
$ws.Range("E29").Value = 0.000777116390577711
$ws.Range("B30").Value = 2959.117122042539
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 0.0295946834879278
$ws.Range("E30").Value = 0.0007853605849685871
$ws.Range("B31").Value = 2982.864196645418
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 0.004497553824902543
$ws.Range("E31").Value = 0.0004959344000781282
$ws.Range("B32").Value = 2998.883965523738
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0.02053416175141775
$ws.Range("E32").Value = 0.004164270564972897
$ws.Range("B33").Value = 3007.118657087335
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 0.01581854507669878
$ws.Range("E33").Value = 0.001106551821162681
$ws.Range("B34").Value = 3064
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 0.02124980265850607
$ws.Range("E34").Value = 0.0008424112900538548
$ws.Range("B35").Value = 3081.81392646662
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0.01832454557534423
$ws.Range("E35").Value = 0.003371716385863408
$ws.Range("B36").Value = 3096.712914411186
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 0.01684630199387963
$ws.Range("E36").Value = 0.00231475905259411
$ws.Range("B37").Value = 3133.618862442508
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 0.005883791576448149
$ws.Range("E37").Value = 0.0004212561682888503
$ws.Range("B38").Value = 3155.201659680821
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 0.008069302435520811
$ws.Range("E38").Value = 0.002241472898755869
$ws.Range("B39").Value = 3168.125784376795
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = 0.00868968257222634
$ws.Range("E39").Value = 0.0009772735551864284
$ws.Range("B40").Value = 3180.855405406395
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 0.0165274208430777
$ws.Range("E40").Value = 0.002253739205874202
$ws.Range("B41").Value = 3192.810075633211
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 0.01624195157086227
$ws.Range("E41").Value = 0.00272765598899971
$ws.Range("B42").Value = 3224.072026396335
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = 0.002810707993906941
$ws.Range("E42").Value = 0.0005110378170739893
$ws.Range("B43").Value = 3243.873355870092
$ws.Range("C43").Value = 2
$ws.Range("D43").Value = 0.002619619068571914
$ws.Range("E43").Value = 0.0003638359817460991
$ws.Range("B44").Value = 3259.461277156505
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 0.007933751544876907
$ws.Range("E44").Value = 0.002319096605425628
$ws.Range("B45").Value = 3275.706563811741
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 0.01189350558226382
$ws.Range("E45").Value = 0.002305883735336952
$ws.Range("B46").Value = 3297.170677130527
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = 0.01137301972618364
$ws.Range("E46").Value = 0.0004465886801115965
$ws.Range("B47").Value = 3322.739257221812
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 0.003987701300459387
$ws.Range("E47").Value = 0.0003063996735747506
$ws.Range("B48").Value = 3340.274605893868
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 0.008777767829112005
$ws.Range("E48").Value = 0.0003701422230533753
$ws.Range("B49").Value = 3363.000857617731
$ws.Range("C49").Value = 3
$ws.Range("D49").Value = 0.01440361255882625
$ws.Range("E49").Value = 0.00318768474662552
$ws.Range("B50").Value = 3380.649591047521
$ws.Range("C50").Value = 3
$ws.Range("D50").Value = 0.01583630578325614
$ws.Range("E50").Value = 0.003401873094180892
$ws.Range("B51").Value = 3395.65849793409
$ws.Range("C51").Value = 3
$ws.Range("D51").Value = 0.006184103894049982
$ws.Range("E51").Value = 0.005250654249665119
$ws.Range("B52").Value = 3411.582461109694
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 0.0130488628990849
$ws.Range("E52").Value = 0.0004610066028577759
$ws.Range("B53").Value = 3430.405841076768
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 0.008578015564787751
$ws.Range("E53").Value = 0.0003750479024411025
$ws.Range("B54").Value = 3452.04083863449
$ws.Range("C54").Value = 3
$ws.Range("D54").Value = 0.01374889714007316
$ws.Range("E54").Value = 0.002062334571010936
$ws.Range("B55").Value = 3477.521886486277
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 0.009769558508139377
$ws.Range("E55").Value = 0.002814957536243618
$ws.Range("B56").Value = 3485.001985506649
$ws.Range("C56").Value = 3
$ws.Range("D56").Value = 0.01542124441400599
$ws.Range("E56").Value = 0.002041047054794886
$ws.Range("B57").Value = 3503.741363695905
$ws.Range("C57").Value = 3
$ws.Range("D57").Value = 0.03074551816829101
$ws.Range("E57").Value = 0.002365039859099297
$ws.Range("B58").Value = 3521.461101658533
$ws.Range("C58").Value = 1
$ws.Range("D58").Value = 0.005014480526134927
$ws.Range("E58").Value = 0.0006412800435417643
$ws.Range("B59").Value = 3536.043535477975
$ws.Range("C59").Value = 1
$ws.Range("D59").Value = 0.003838362625454867
$ws.Range("E59").Value = 0.0006414377843615218
$ws.Range("B60").Value = 3545.218978961481
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = 0.01910140741942075
$ws.Range("E60").Value = 0.0007450977178528679
$ws.Range("B61").Value = 3561.667933967138
$ws.Range("C61").Value = 3
$ws.Range("D61").Value = 0.01722724925451776
$ws.Range("E61").Value = 0.002429483869226814
$ws.Range("B62").Value = 3580.019183465469
$ws.Range("C62").Value = 3
$ws.Range("D62").Value = 0.008285624013533003
$ws.Range("E62").Value = 0.0006833320461825788
$ws.Range("B63").Value = 3589.932337372159
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 0.04782351119169598
$ws.Range("E63").Value = 0.00452149560357852

Write-Output "done"